# Generate Report for Handback
#
# Marks the two localized files (zh-cn, de-de) as handed back: updates the
# Status column text, stamps the "Latest Target File" / "Latest Handback
# File" hyperlinked columns (F/G) with the file names that were handed
# back, and records the real handback datetime in column H (replacing the
# "0001-01-01 00:00:00" placeholder).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: just the status text changes (cols B and C) -----------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# --- Per-locale detail sheets ----------------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn"
        HandbackDate = "2016-03-20 20:39:03"
        Row2Md = "3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.md"
        Row2MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/355c163127c7faf33e52d1c6c08328993d671750/e2e/3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.md"
        Row2Xlf = "3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.ea38419ac826560141e29a7cfbdb93fb872834d8.zh-cn.xlf"
        Row2XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d846a3c0847067811e0dce1989f4c0f80f2d876/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.ea38419ac826560141e29a7cfbdb93fb872834d8.zh-cn.xlf"
        Row3Md = "502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.md"
        Row3MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/355c163127c7faf33e52d1c6c08328993d671750/e2e/502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.md"
        Row3Xlf = "502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.48e9be9db298239a2ae54e57d5c2d8184d09ff4d.zh-cn.xlf"
        Row3XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d846a3c0847067811e0dce1989f4c0f80f2d876/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.48e9be9db298239a2ae54e57d5c2d8184d09ff4d.zh-cn.xlf"
    },
    @{
        Sheet = "de-de"
        HandbackDate = "2016-03-20 20:39:10"
        Row2Md = "3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.md"
        Row2MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/355c163127c7faf33e52d1c6c08328993d671750/e2e/3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.md"
        Row2Xlf = "3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.ea38419ac826560141e29a7cfbdb93fb872834d8.de-de.xlf"
        Row2XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4edbc844406c55a7d6ccb55c9d6e7d71599b7f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3d9a67f5-8734-414c-bcc2-2d6fff5a0a99.ea38419ac826560141e29a7cfbdb93fb872834d8.de-de.xlf"
        Row3Md = "502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.md"
        Row3MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/355c163127c7faf33e52d1c6c08328993d671750/e2e/502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.md"
        Row3Xlf = "502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.48e9be9db298239a2ae54e57d5c2d8184d09ff4d.de-de.xlf"
        Row3XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4edbc844406c55a7d6ccb55c9d6e7d71599b7f8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/502a4e9c-e0f6-48d4-bfa2-6d6278c5b9bc.48e9be9db298239a2ae54e57d5c2d8184d09ff4d.de-de.xlf"
    }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # New "Latest Target File" (F) / "Latest Handback File" (G) hyperlinked
    # columns for row 2 (the 3d9a67f5... source file).
    $ws.Hyperlinks.Add($ws.Range("F2"), $loc.Row2MdUrl, "", "", $loc.Row2Md)
    $ws.Hyperlinks.Add($ws.Range("G2"), $loc.Row2XlfUrl, "", "", $loc.Row2Xlf)

    # ... and for row 3 (the 502a4e9c... source file).
    $ws.Hyperlinks.Add($ws.Range("F3"), $loc.Row3MdUrl, "", "", $loc.Row3Md)
    $ws.Hyperlinks.Add($ws.Range("G3"), $loc.Row3XlfUrl, "", "", $loc.Row3Xlf)

    # Latest Handback DateTime (H), replacing the "0001-01-01 00:00:00"
    # placeholder with the real handback timestamp.
    $ws.Range("H2").Value = $loc.HandbackDate
    $ws.Range("H3").Value = $loc.HandbackDate
}
